# Update workbook: refresh COVID-19 country data and fix country row order
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp in the title row (A1)
$ws.Range("A1").Value = "Datos actualizados a 28 de Mayo de 2020 a las 23:10"

# Row 4: 'Estados Unidos' -> 'Estados Unidos'
$ws.Range("B4").Value = 1765057
$ws.Range("C4").Value = 19254
$ws.Range("D4").Value = 497087
$ws.Range("E4").Value = 1164764
$ws.Range("G4").Value = 1099
$ws.Range("H4").Value = 103206

# Row 16: 'Canada' -> 'Canada'
$ws.Range("B16").Value = 88475
$ws.Range("C16").Value = 956
$ws.Range("D16").Value = 46777
$ws.Range("E16").Value = 34825

# Row 106: 'Guinea-Bisau' -> 'Guinea-Bisau'
$ws.Range("E106").Value = 1145
$ws.Range("G106").Value = 1
$ws.Range("H106").Value = 8

# Row 118: 'Niger' -> 'Niger'
$ws.Range("D118").Value = 803
$ws.Range("E118").Value = 88

# Row 140: 'Estado de Palestina' -> 'Estado de Palestina'
$ws.Range("D140").Value = 368
$ws.Range("E140").Value = 75

# Row 172: 'Aruba' -> 'Libia'
$ws.Range("A172").Value = "Libia"
$ws.Range("B172").Value = 105
$ws.Range("C172").Value = 6
$ws.Range("D172").Value = 40
$ws.Range("E172").Value = 60
$ws.Range("G172").Value = 1
$ws.Range("H172").Value = 5

# Row 173: 'Bahamas' -> 'Aruba'
$ws.Range("A173").Value = "Aruba"
$ws.Range("B173").Value = 101
$ws.Range("D173").Value = 97
$ws.Range("E173").Value = 1
$ws.Range("H173").Value = 3

# Row 174: 'Libia' -> 'Bahamas'
$ws.Range("A174").Value = "Bahamas"
$ws.Range("B174").Value = 100
$ws.Range("D174").Value = 46
$ws.Range("E174").Value = 43
$ws.Range("H174").Value = 11

# Row 180: 'Angola' -> 'Angola'
$ws.Range("B180").Value = 74
$ws.Range("C180").Value = 3
$ws.Range("E180").Value = 52

# Row 198: 'Curazao' -> 'Fiyi'
$ws.Range("A198").Value = "Fiyi"
$ws.Range("D198").Value = 15
$ws.Range("H198").Value = 0

# Row 199: 'Fiyi' -> 'Curazao'
$ws.Range("A199").Value = "Curazao"
$ws.Range("D199").Value = 14
$ws.Range("H199").Value = 1

# Row 200: 'Santa Lucia' -> 'Belice'
$ws.Range("A200").Value = "Belice"
$ws.Range("D200").Value = 16
$ws.Range("H200").Value = 2

# Row 201: 'Belice' -> 'Santa Lucia'
$ws.Range("A201").Value = "Santa Lucia"
$ws.Range("D201").Value = 18
$ws.Range("H201").Value = 0

# Row 210: 'Seychelles' -> 'Montserrat'
$ws.Range("A210").Value = "Montserrat"
$ws.Range("D210").Value = 10
$ws.Range("H210").Value = 1

# Row 211: 'Montserrat' -> 'Seychelles'
$ws.Range("A211").Value = "Seychelles"
$ws.Range("D211").Value = 11
$ws.Range("H211").Value = 0

# Row 215: 'Bonaire, San Eustaquio y Saba' -> 'San Bartolome'
$ws.Range("A215").Value = "San Bartolome"

# Row 216: 'San Bartolome' -> 'Bonaire, San Eustaquio y Saba'
$ws.Range("A216").Value = "Bonaire, San Eustaquio y Saba"

